$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$noteRange = $ws1.Range("A1")
$note = $noteRange.Value2
$note = $note.Replace("1000 Bs = 7.04 = 27756.62 pesos", "1000 Bs = 6.96 = 27582.61 pesos")
$note = $note.Replace("27756.62 pesos = 7.03 = 948.64 Bs", "27582.61 pesos = 6.95 = 972.69 Bs")
$noteRange.Value = $note

# --- Update the transfi rate table on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 143.75
$ws2.Range("O10").Value = 3965
$ws2.Range("N12").Value = 3970
$ws2.Range("O12").Value = 140
